# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> currently the stock "Office Theme" (clrScheme "Office")
#   ppt/theme/theme2.xml  -> currently the "Integral" theme (clrScheme "Red Violet"),
#                             used by the slide master / the whole visible deck
#
# The authored change swaps the two themes' contents: the deck's live theme
# (theme2.xml, reached through SlideMaster.Theme / Slide.ThemeColorScheme)
# becomes the plain "Office Theme" palette, and theme1.xml becomes the
# "Integral"/Red-Violet palette. The font scheme and the effect/format scheme
# are identical between the two themes already, so the only observable
# difference is the twelve theme colors (dk1/lt1/dk2/lt2/accent1-6/hlink/
# folHlink) plus the theme/colour-scheme display names.
#
# Apply the reachable half of that swap: push the standard Office Theme
# colors onto the presentation's live color scheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink, in that fixed order).

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0          # dk1      -> 000000
$cs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$cs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  -> FFC000
$cs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$cs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$cs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$cs.Item(12).RGB = 7491477    # folHlink -> 954F72
